$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.686.56'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '1.893.00'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'238.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = "'0.4839"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('D8').Value = "'0.2886"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.02%  '
$ws.Range('D9').Value = "'0.06551"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('D10').Value = '1.838.80'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('D12').Value = "'0.07463"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('D13').Value = "'5.112"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').Value = "'88.01"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = "'0.6682"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '30.645.50'
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('D17').Value = "'13.25"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = "'0.000007594"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = "'232.94"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.80%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.284"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.069.98'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = "'6.194"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').Value = "'9.397"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range('D26').Value = "'168.94"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.26%  '
$ws.Range('D27').Value = "'18.74"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.75%  '
$ws.Range('D28').Value = "'1.958"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.92%  '
$ws.Range('D29').Value = "'0.1024"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.57%  '
$ws.Range('D30').Value = "'1.399"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('D31').Value = "'4.338"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('D32').Value = "'4.028"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  +2.17%  '
$ws.Range('E34').Value = '  +5.92%  '
$ws.Range('E35').Value = '  +4.05%  '
$ws.Range('D36').Value = "'0.9991"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = "'2.715"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('D38').Value = "'0.01886"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.64%  '
$ws.Range('D39').Value = "'2.650"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('D40').Value = "'0.9204"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').Value = "'2.070"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.91%  '
$ws.Range('D42').Value = "'106.83"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.29%  '
$ws.Range('D43').Value = "'0.4291"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('D45').Value = "'5.633"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.96%  '
$ws.Range('D46').Value = "'7.426"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.40%  '
$ws.Range('D47').Value = "'64.21"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('D49').Value = "'1.490"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Value = "'9.010"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('E51').Value = '  +1.35%  '
